$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "6-5-5"
$ws.Range("C2").Value = "new_sequential"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.166"
$ws.Range("E2").Value = "sat"
$ws.Range("F2").Value = 6276
$ws.Range("G2").Value = 170692
